$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.169
$ws.Range("E3").Value = 16.856
$ws.Range("A12").Value = -21.651
$ws.Range("D14").Value = -7.716000000000001
$ws.Range("D26").Value = -7.673
$ws.Range("E30").Value = 16.323
$ws.Range("D31").Value = -7.915000000000002
$ws.Range("A32").Value = -21.761
$ws.Range("D35").Value = -7.939
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.741
$ws.Range("A38").Value = -19.741
$ws.Range("E44").Value = 16.611
$ws.Range("D45").Value = -7.539
$ws.Range("A46").Value = -21.849
$ws.Range("A54").Value = -21.869
$ws.Range("A55").Value = -22.135
$ws.Range("D57").Value = -8.271000000000001
$ws.Range("E58").Value = 16.65
$ws.Range("A67").Value = -21.577
$ws.Range("A69").Value = -21.72
$ws.Range("A72").Value = -21.445
$ws.Range("E84").Value = 16.415
$ws.Range("E89").Value = 17.149
$ws.Range("A91").Value = -21.586
$ws.Range("E91").Value = 17.243
$ws.Range("E92").Value = 17.052
$ws.Range("A99").Value = -20.43
$ws.Range("D100").Value = -8.280000000000001
$ws.Range("D102").Value = -7.606
$ws.Range("E102").Value = 16.552
